$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = "A1BG"
$ws.Cells.Item(22, 2).Value = "ENSG00000121410.11"
$ws.Cells.Item(22, 3).Value = 6.45
$ws.Cells.Item(22, 4).Value = 18.61
$ws.Cells.Item(22, 5).Value = -1.396
$ws.Cells.Item(22, 6).Value = 0.00000000000000000000000000000401

$ws.Cells.Item(23, 1).Value = "AACSP1"
$ws.Cells.Item(23, 2).Value = "ENSG00000250420.8"
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 2.535
$ws.Cells.Item(23, 5).Value = -1.822
$ws.Cells.Item(23, 6).Value = [double]"4.1e-229"

$ws.Cells.Item(24, 1).Value = "AAGAB"
$ws.Cells.Item(24, 2).Value = "ENSG00000103591.12"
$ws.Cells.Item(24, 3).Value = 20.52
$ws.Cells.Item(24, 4).Value = 52.548
$ws.Cells.Item(24, 5).Value = -1.315
$ws.Cells.Item(24, 6).Value = 0.00000000000000000000000000000000000000000000000000000000000000398

$ws.Cells.Item(25, 1).Value = "AAMDC"
$ws.Cells.Item(25, 2).Value = "ENSG00000087884.14"
$ws.Cells.Item(25, 3).Value = 13.39
$ws.Cells.Item(25, 4).Value = 28.505
$ws.Cells.Item(25, 5).Value = -1.036
$ws.Cells.Item(25, 6).Value = 0.000000000000000000000000000000000000000000000000366

$ws.Cells.Item(26, 1).Value = "AARS"
$ws.Cells.Item(26, 2).Value = "ENSG00000090861.15"
$ws.Cells.Item(26, 3).Value = 31.441
$ws.Cells.Item(26, 4).Value = 205.329
$ws.Cells.Item(26, 5).Value = -2.669
$ws.Cells.Item(26, 6).Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000238

$ws.Cells.Item(27, 1).Value = "ABALON"
$ws.Cells.Item(27, 2).Value = "ENSG00000281376.1"
$ws.Cells.Item(27, 3).Value = 0.96
$ws.Cells.Item(27, 4).Value = 3.3
$ws.Cells.Item(27, 5).Value = -1.133
$ws.Cells.Item(27, 6).Value = 0.000000000000000000000000000000000000000000000000000000000000128

$ws.Cells.Item(28, 1).Value = "ABCA3"
$ws.Cells.Item(28, 2).Value = "ENSG00000167972.13"
$ws.Cells.Item(28, 3).Value = 2.38
$ws.Cells.Item(28, 4).Value = 7.385
$ws.Cells.Item(28, 5).Value = -1.311
$ws.Cells.Item(28, 6).Value = 0.00000000000000000000000264

$ws.Cells.Item(29, 1).Value = "ABCA8"
$ws.Cells.Item(29, 2).Value = "ENSG00000141338.13"
$ws.Cells.Item(29, 3).Value = 0.16
$ws.Cells.Item(29, 4).Value = 12.8
$ws.Cells.Item(29, 5).Value = -3.572
$ws.Cells.Item(29, 6).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000231

$ws.Cells.Item(30, 1).Value = "ABCB10P1"
$ws.Cells.Item(30, 2).Value = "ENSG00000274099.1"
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 2.75
$ws.Cells.Item(30, 5).Value = -1.907
$ws.Cells.Item(30, 6).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000012

$ws.Cells.Item(31, 1).Value = "ABCB10P3"
$ws.Cells.Item(31, 2).Value = "ENSG00000261524.1"
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 10.505
$ws.Cells.Item(31, 5).Value = -3.524
$ws.Cells.Item(31, 6).Value = [double]"4.13e-251"

$ws.Cells.Item(32, 1).Value = "ABCB10P4"
$ws.Cells.Item(32, 2).Value = "ENSG00000260053.2"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 25.325
$ws.Cells.Item(32, 5).Value = -4.718
$ws.Cells.Item(32, 6).Value = [double]"1.42e-256"

$ws.Cells.Item(33, 1).Value = "ABCB6"
$ws.Cells.Item(33, 2).Value = "ENSG00000115657.12"
$ws.Cells.Item(33, 3).Value = 5.31
$ws.Cells.Item(33, 4).Value = 20.2
$ws.Cells.Item(33, 5).Value = -1.748
$ws.Cells.Item(33, 6).Value = 0.00000000000000000000000000000134

$ws.Cells.Item(34, 1).Value = "ABCB8"
$ws.Cells.Item(34, 2).Value = "ENSG00000197150.12"
$ws.Cells.Item(34, 3).Value = 18.02
$ws.Cells.Item(34, 4).Value = 45.984
$ws.Cells.Item(34, 5).Value = -1.305
$ws.Cells.Item(34, 6).Value = 0.00000000000000000000000000000000000000000000000000021

$ws.Cells.Item(35, 1).Value = "ABCC4"
$ws.Cells.Item(35, 2).Value = "ENSG00000125257.13"
$ws.Cells.Item(35, 3).Value = 9.07
$ws.Cells.Item(35, 4).Value = 20.554
$ws.Cells.Item(35, 5).Value = -1.098
$ws.Cells.Item(35, 6).Value = 0.000000000000000246

$ws.Cells.Item(36, 1).Value = "ABCF2"
$ws.Cells.Item(36, 2).Value = "ENSG00000033050.7"
$ws.Cells.Item(36, 3).Value = 24.609
$ws.Cells.Item(36, 4).Value = 54.57
$ws.Cells.Item(36, 5).Value = -1.118
$ws.Cells.Item(36, 6).Value = 0.000000000000000000000000000000000000000000000000000063

$ws.Cells.Item(37, 1).Value = "ABHD11"
$ws.Cells.Item(37, 2).Value = "ENSG00000106077.18"
$ws.Cells.Item(37, 3).Value = 7.67
$ws.Cells.Item(37, 4).Value = 21.41
$ws.Cells.Item(37, 5).Value = -1.37
$ws.Cells.Item(37, 6).Value = 0.00000000000000000000000000000000000545

$ws.Cells.Item(38, 1).Value = "ABHD5"
$ws.Cells.Item(38, 2).Value = "ENSG00000011198.7"
$ws.Cells.Item(38, 3).Value = 13.85
$ws.Cells.Item(38, 4).Value = 38.019
$ws.Cells.Item(38, 5).Value = -1.394
$ws.Cells.Item(38, 6).Value = 0.000000000000000000000000000000000104

$ws.Cells.Item(39, 1).Value = "ABO"
$ws.Cells.Item(39, 2).Value = "ENSG00000175164.13"
$ws.Cells.Item(39, 3).Value = 1.52
$ws.Cells.Item(39, 4).Value = 10.93
$ws.Cells.Item(39, 5).Value = -2.243
$ws.Cells.Item(39, 6).Value = 0.0000000000103

$ws.Cells.Item(40, 1).Value = "ABT1"
$ws.Cells.Item(40, 2).Value = "ENSG00000146109.4"
$ws.Cells.Item(40, 3).Value = 9.66
$ws.Cells.Item(40, 4).Value = 24.48
$ws.Cells.Item(40, 5).Value = -1.257
$ws.Cells.Item(40, 6).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000217

$ws.Cells.Item(41, 1).Value = "AC000041.8"
$ws.Cells.Item(41, 2).Value = "ENSG00000242156.1"
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 7.885
$ws.Cells.Item(41, 5).Value = -3.151
$ws.Cells.Item(41, 6).Value = 0.000000000000000000000000000000135
